# Work-Time.xlsx: rename the "Zeit" column header to "Zeit in h" and log the
# first entry's worked hours (2h) for 2021-08-16, applying the matching date
# number format to the date column/cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header B1: "Zeit" -> "Zeit in h"
$ws.Range("B1").Value = "Zeit in h"

# 2. Log 2 hours worked in B2 (row for the 2021-08-16 entry in A2)
$ws.Range("B2").Value = 2

# 3. Give the date column its date formatting. A2 is formatted first so its
#    existing style slot is reused, then A1 (the header cell) picks up the
#    same number format while keeping its own border/style.
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A1").NumberFormat = "mm-dd-yy"
